# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Black Amber" plums at the top of the
# existing block (new rows 117-119), shifting the previous rows 117-142
# down to 120-145.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 117.
$ws.Range("A117:T119").EntireRow.Insert()

# Common column values shared by every data row in this block.
$mercadoId   = 5
$mercado     = 'Macroferia Regional de Talca'
$region      = 'Maule'
$codreg      = 7
$tipo        = 'Fruta'
$productoId  = 100103
$producto    = 'Frutos de hueso (carozo)'
$categoriaId = 100103002
$categoria   = 'Ciruela'
$unidad      = '$/bandeja 18 kilos granel'
$origen      = 'Provincia de Curicó'
$kgUnidad    = 18

# New row 117: Black Amber / Especial
$r = 117
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = 44946
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value2  = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = 'Black Amber'
$ws.Cells.Item($r, 12).Value2 = 'Especial'
$ws.Cells.Item($r, 13).Value2 = 150
$ws.Cells.Item($r, 14).Value2 = 12000
$ws.Cells.Item($r, 15).Value2 = 12000
$ws.Cells.Item($r, 16).Value2 = 12000
$ws.Cells.Item($r, 17).Value2 = $unidad
$ws.Cells.Item($r, 18).Value2 = $origen
$ws.Cells.Item($r, 19).Value2 = 667
$ws.Cells.Item($r, 20).Value2 = $kgUnidad

# New row 118: Black Amber / Primera
$r = 118
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = 44946
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value2  = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = 'Black Amber'
$ws.Cells.Item($r, 12).Value2 = 'Primera'
$ws.Cells.Item($r, 13).Value2 = 250
$ws.Cells.Item($r, 14).Value2 = 10000
$ws.Cells.Item($r, 15).Value2 = 10000
$ws.Cells.Item($r, 16).Value2 = 10000
$ws.Cells.Item($r, 17).Value2 = $unidad
$ws.Cells.Item($r, 18).Value2 = $origen
$ws.Cells.Item($r, 19).Value2 = 556
$ws.Cells.Item($r, 20).Value2 = $kgUnidad

# New row 119: Black Amber / Segunda
$r = 119
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = 44946
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value2  = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = 'Black Amber'
$ws.Cells.Item($r, 12).Value2 = 'Segunda'
$ws.Cells.Item($r, 13).Value2 = 200
$ws.Cells.Item($r, 14).Value2 = 8000
$ws.Cells.Item($r, 15).Value2 = 8000
$ws.Cells.Item($r, 16).Value2 = 8000
$ws.Cells.Item($r, 17).Value2 = $unidad
$ws.Cells.Item($r, 18).Value2 = $origen
$ws.Cells.Item($r, 19).Value2 = 444
$ws.Cells.Item($r, 20).Value2 = $kgUnidad

$ws.Range('D117:D119').NumberFormat = 'YYYY-MM-DD HH:MM:SS'
